# Trade #40 closed at 2026-02-17 12:47:25 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.04
$wsSummary.Range("B4").Value = 0.03
$wsSummary.Range("B5").Value = 0.01
$wsSummary.Range("B6").Value = 40
$wsSummary.Range("B8").Value = 16
$wsSummary.Range("B9").Value = 37.5

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.04
$wsStatus.Range("D4").Value = 40
$wsStatus.Range("E4").Value = 0.03
$wsStatus.Range("F4").Value = 0.04
$wsStatus.Range("G4").Value = 37.5

# --- All Trades sheet (Trade #40 -> row 41) ---
$wsTrades = $wb.Worksheets.Item("All Trades")
$wsTrades.Range("G41").Value = 0.2
$wsTrades.Range("H41").Value = "CLOSED"
$wsTrades.Range("I41").Value = -56.5217
$wsTrades.Range("J41").Value = -0.26
$wsTrades.Range("K41").Value = 100.04
$wsTrades.Range("P41").Value = "early_exit"
$wsTrades.Range("Q41").Value = 5.53

# --- MarketMaking sheet (same Trade #40 -> row 41) ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G41").Value = 0.2
$wsMM.Range("H41").Value = "CLOSED"
$wsMM.Range("I41").Value = -56.5217
$wsMM.Range("J41").Value = -0.26
$wsMM.Range("K41").Value = 100.04
$wsMM.Range("P41").Value = "early_exit"
$wsMM.Range("Q41").Value = 5.53
